$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (45) with the next month's data (01-09-2021)
# Force the cell to text format so Excel does not auto-convert the
# dd-mm-yyyy-looking string into a date serial number, then restore the
# cell's style to the default "Normal" style so it keeps using the
# workbook's default (unstyled) formatting, matching the other date
# cells in column A (A2:A44).
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "01-09-2021"
$ws.Range("A45").Style = "Normal"
$ws.Range("B45").Value = 1.2
$ws.Range("C45").Value = 0.9
$ws.Range("D45").Value = 1.5
$ws.Range("E45").Value = 0.8
$ws.Range("F45").Value = 3.9
